$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

Set-TextValue 'D2' '29.395.51'
$ws.Range('E2').Value = '  -0.12%  '
Set-TextValue 'D3' '1.845.41'
$ws.Range('E3').Value = '  -0.25%  '
Set-TextValue 'D4' '0.9994'
$ws.Range('E4').Value = '  +0.03%  '
Set-TextValue 'D5' '239.12'
$ws.Range('E5').Value = '  -0.79%  '
Set-TextValue 'D6' '0.6322'
$ws.Range('E6').Value = '  +0.03%  '
Set-TextValue 'D7' '1.001'
$ws.Range('E7').Value = '  +0.04%  '
Set-TextValue 'D8' '0.07571'
$ws.Range('E8').Value = '  +0.10%  '
$ws.Range('E9').Value = '  -1.00%  '
$ws.Range('E10').Value = '  -0.05%  '
Set-TextValue 'D11' '0.07716'
$ws.Range('E11').Value = '  -0.13%  '
Set-TextValue 'D12' '1.852.42'
$ws.Range('E12').Value = '  -6.68%  '
Set-TextValue 'D13' '5.005'
$ws.Range('E13').Value = '  +0.29%  '
Set-TextValue 'D14' '0.6801'
$ws.Range('E14').Value = '  -0.55%  '
$ws.Range('E15').Value = '  +5.89%  '
Set-TextValue 'D16' '83.36'
$ws.Range('E16').Value = '  +0.49%  '
Set-TextValue 'D17' '2.086.61'
$ws.Range('E17').Value = '  -7.85%  '
Set-TextValue 'D18' '6.177'
$ws.Range('E18').Value = '  -0.23%  '
Set-TextValue 'D19' '29.428.56'
$ws.Range('E19').Value = '  -0.14%  '
Set-TextValue 'D20' '229.01'
$ws.Range('E20').Value = '  -1.06%  '
$ws.Range('E21').Value = '  -0.37%  '
Set-TextValue 'D22' '1.000'
$ws.Range('E22').Value = '  +0.03%  '
Set-TextValue 'D23' '7.470'
$ws.Range('E23').Value = '  -1.64%  '
Set-TextValue 'D24' '1.001'
$ws.Range('E24').Value = '  +0.05%  '
Set-TextValue 'D25' '156.94'
$ws.Range('E25').Value = '  +0.72%  '
$ws.Range('E26').Value = '  +0.83%  '
Set-TextValue 'D27' '8.364'
$ws.Range('E27').Value = '  -0.45%  '
Set-TextValue 'D28' '17.59'
$ws.Range('E28').Value = '  -0.56%  '
$ws.Range('B29').Value = 'PancakeSwap'
$ws.Range('C29').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
Set-TextValue 'D29' '1.459'
$ws.Range('E29').Value = '  -0.86%  '
$ws.Range('B30').Value = 'Toncoin'
$ws.Range('C30').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-TextValue 'D30' '1.298'
$ws.Range('E30').Value = '  +3.19%  '
Set-TextValue 'D31' '0.05648'
$ws.Range('E31').Value = '  -2.31%  '
Set-TextValue 'D32' '4.101'
$ws.Range('E32').Value = '  -0.71%  '
Set-TextValue 'D33' '4.023'
$ws.Range('E33').Value = '  +0.11%  '
Set-TextValue 'D34' '1.847'
$ws.Range('E34').Value = '  -0.27%  '
$ws.Range('E35').Value = '  -0.14%  '
Set-TextValue 'D36' '0.7102'
$ws.Range('E36').Value = '  -0.89%  '
Set-TextValue 'D37' '2.592'
$ws.Range('E37').Value = '  -0.18%  '
Set-TextValue 'D38' '1.248.58'
$ws.Range('E38').Value = '  -0.18%  '
Set-TextValue 'D39' '0.01810'
$ws.Range('E39').Value = '  +0.26%  '
Set-TextValue 'D40' '2.769'
$ws.Range('E40').Value = '  -1.14%  '
Set-TextValue 'D41' '6.383'
$ws.Range('E41').Value = '  +4.77%  '
Set-TextValue 'D42' '0.9019'
$ws.Range('E42').Value = '  -0.28%  '
$ws.Range('E43').Value = '  +0.05%  '
Set-TextValue 'D44' '101.79'
$ws.Range('E44').Value = '  +0.09%  '
Set-TextValue 'D45' '65.84'
$ws.Range('E45').Value = '  -1.64%  '
$ws.Range('B46').Value = 'BabyDogeCoin'
$ws.Range('C46').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
Set-TextValue 'D46' '0.00000000119'
$ws.Range('E46').Value = '  +0.46%  '
$ws.Range('B47').Value = 'Aptos'
$ws.Range('C47').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextValue 'D47' '7.098'
$ws.Range('E47').Value = '  -0.78%  '
$ws.Range('E48').Value = '  -0.48%  '
$ws.Range('E49').Value = '  -0.40%  '
Set-TextValue 'D50' '8.924'
$ws.Range('E50').Value = '  -2.94%  '
$ws.Range('E51').Value = '  -0.13%  '
